$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 69.89967366666667
$ws.Range("H2").Value = 209.699021
$ws.Range("I2").Value = 0.6608367681537789
$ws.Range("J2").Value = 0.660836768153779
$ws.Range("M2").Value = 45.90594266666667
$ws.Range("N2").Value = 137.717828
$ws.Range("O2").Value = 0.3954672001633582
$ws.Range("P2").Value = 0.3954672001633583
$ws.Range("Q2").Value = 3208.81041176071
$ws.Range("R2").Value = 28879.29370584639
$ws.Range("S2").Value = 0.2613392664667772
$ws.Range("T2").Value = 0.2613392664667773

$ws.Range("G3").Value = 69.89967366666667
$ws.Range("H3").Value = 209.699021
$ws.Range("I3").Value = 0.6608367681537789
$ws.Range("J3").Value = 0.660836768153779
$ws.Range("M3").Value = 40.44578266666667
$ws.Range("N3").Value = 121.337348
$ws.Range("O3").Value = 0.3484294080560655
$ws.Range("P3").Value = 0.3484294080560656
$ws.Range("Q3").Value = 2827.147009592923
$ws.Range("R3").Value = 25444.32308633631
$ws.Range("S3").Value = 0.2302549639495046
$ws.Range("T3").Value = 0.2302549639495047

$ws.Range("G4").Value = 69.89967366666667
$ws.Range("H4").Value = 209.699021
$ws.Range("I4").Value = 0.6608367681537789
$ws.Range("J4").Value = 0.660836768153779
$ws.Range("M4").Value = 12.761795
$ws.Range("N4").Value = 38.28538500000001
$ws.Range("O4").Value = 0.1099393900775594
$ws.Range("P4").Value = 0.1099393900775594
$ws.Range("Q4").Value = 892.0453059008985
$ws.Range("R4").Value = 8028.407753108087
$ws.Range("S4").Value = 0.07265199123165199
$ws.Range("T4").Value = 0.07265199123165202

$ws.Range("G5").Value = 69.89967366666667
$ws.Range("H5").Value = 209.699021
$ws.Range("I5").Value = 0.6608367681537789
$ws.Range("J5").Value = 0.660836768153779
$ws.Range("M5").Value = 16.966758
$ws.Range("N5").Value = 50.900274
$ws.Range("O5").Value = 0.1461640017030168
$ws.Range("P5").Value = 0.1461640017030168
$ws.Range("Q5").Value = 1185.970847381306
$ws.Range("R5").Value = 10673.73762643175
$ws.Range("S5").Value = 0.09659054650584506
$ws.Range("T5").Value = 0.09659054650584509

$ws.Range("G6").Value = 17.57434666666667
$ws.Range("H6").Value = 52.72304
$ws.Range("I6").Value = 0.1661491941864736
$ws.Range("J6").Value = 0.1661491941864736
$ws.Range("M6").Value = 45.90594266666667
$ws.Range("N6").Value = 137.717828
$ws.Range("O6").Value = 0.3954672001633582
$ws.Range("P6").Value = 0.3954672001633583
$ws.Range("Q6").Value = 806.7669504841244
$ws.Range("R6").Value = 7260.90255435712
$ws.Range("S6").Value = 0.06570655663432283
$ws.Range("T6").Value = 0.06570655663432286

$ws.Range("G7").Value = 17.57434666666667
$ws.Range("H7").Value = 52.72304
$ws.Range("I7").Value = 0.1661491941864736
$ws.Range("J7").Value = 0.1661491941864736
$ws.Range("M7").Value = 40.44578266666667
$ws.Range("N7").Value = 121.337348
$ws.Range("O7").Value = 0.3484294080560655
$ws.Range("P7").Value = 0.3484294080560656
$ws.Range("Q7").Value = 710.8082057886578
$ws.Range("R7").Value = 6397.273852097919
$ws.Range("S7").Value = 0.05789126537938528
$ws.Range("T7").Value = 0.05789126537938531

$ws.Range("G8").Value = 17.57434666666667
$ws.Range("H8").Value = 52.72304
$ws.Range("I8").Value = 0.1661491941864736
$ws.Range("J8").Value = 0.1661491941864736
$ws.Range("M8").Value = 12.761795
$ws.Range("N8").Value = 38.28538500000001
$ws.Range("O8").Value = 0.1099393900775594
$ws.Range("P8").Value = 0.1099393900775594
$ws.Range("Q8").Value = 224.2802094189334
$ws.Range("R8").Value = 2018.5218847704
$ws.Range("S8").Value = 0.01826634107073889
$ws.Range("T8").Value = 0.0182663410707389

$ws.Range("G9").Value = 17.57434666666667
$ws.Range("H9").Value = 52.72304
$ws.Range("I9").Value = 0.1661491941864736
$ws.Range("J9").Value = 0.1661491941864736
$ws.Range("M9").Value = 16.966758
$ws.Range("N9").Value = 50.900274
$ws.Range("O9").Value = 0.1461640017030168
$ws.Range("P9").Value = 0.1461640017030168
$ws.Range("Q9").Value = 298.17968690144
$ws.Range("R9").Value = 2683.61718211296
$ws.Range("S9").Value = 0.0242850311020266
$ws.Range("T9").Value = 0.02428503110202661

$ws.Range("G10").Value = 4.152730666666667
$ws.Range("H10").Value = 12.458192
$ws.Range("I10").Value = 0.0392602278210887
$ws.Range("J10").Value = 0.03926022782108871
$ws.Range("M10").Value = 45.90594266666667
$ws.Range("N10").Value = 137.717828
$ws.Range("O10").Value = 0.3954672001633582
$ws.Range("P10").Value = 0.3954672001633583
$ws.Range("Q10").Value = 190.6350158941084
$ws.Range("R10").Value = 1715.715143046976
$ws.Range("S10").Value = 0.01552613237418153
$ws.Range("T10").Value = 0.01552613237418153

$ws.Range("G11").Value = 4.152730666666667
$ws.Range("H11").Value = 12.458192
$ws.Range("I11").Value = 0.0392602278210887
$ws.Range("J11").Value = 0.03926022782108871
$ws.Range("M11").Value = 40.44578266666667
$ws.Range("N11").Value = 121.337348
$ws.Range("O11").Value = 0.3484294080560655
$ws.Range("P11").Value = 0.3484294080560656
$ws.Range("Q11").Value = 167.9604420172018
$ws.Range("R11").Value = 1511.643978154816
$ws.Range("S11").Value = 0.01367941793984821
$ws.Range("T11").Value = 0.01367941793984822

$ws.Range("G12").Value = 4.152730666666667
$ws.Range("H12").Value = 12.458192
$ws.Range("I12").Value = 0.0392602278210887
$ws.Range("J12").Value = 0.03926022782108871
$ws.Range("M12").Value = 12.761795
$ws.Range("N12").Value = 38.28538500000001
$ws.Range("O12").Value = 0.1099393900775594
$ws.Range("P12").Value = 0.1099393900775594
$ws.Range("Q12").Value = 52.99629745821334
$ws.Range("R12").Value = 476.9666771239201
$ws.Range("S12").Value = 0.004316245500956521
$ws.Range("T12").Value = 0.004316245500956523

$ws.Range("G13").Value = 4.152730666666667
$ws.Range("H13").Value = 12.458192
$ws.Range("I13").Value = 0.0392602278210887
$ws.Range("J13").Value = 0.03926022782108871
$ws.Range("M13").Value = 16.966758
$ws.Range("N13").Value = 50.900274
$ws.Range("O13").Value = 0.1461640017030168
$ws.Range("P13").Value = 0.1461640017030168
$ws.Range("Q13").Value = 70.458376260512
$ws.Range("R13").Value = 634.1253863446079
$ws.Range("S13").Value = 0.005738432006102436
$ws.Range("T13").Value = 0.005738432006102438

$ws.Range("G14").Value = 14.14774133333333
$ws.Range("H14").Value = 42.443224
$ws.Range("I14").Value = 0.1337538098386587
$ws.Range("J14").Value = 0.1337538098386588
$ws.Range("M14").Value = 45.90594266666667
$ws.Range("N14").Value = 137.717828
$ws.Range("O14").Value = 0.3954672001633582
$ws.Range("P14").Value = 0.3954672001633583
$ws.Range("Q14").Value = 649.4654025108302
$ws.Range("R14").Value = 5845.188622597472
$ws.Range("S14").Value = 0.0528952446880766
$ws.Range("T14").Value = 0.05289524468807662

$ws.Range("G15").Value = 14.14774133333333
$ws.Range("H15").Value = 42.443224
$ws.Range("I15").Value = 0.1337538098386587
$ws.Range("J15").Value = 0.1337538098386588
$ws.Range("M15").Value = 40.44578266666667
$ws.Range("N15").Value = 121.337348
$ws.Range("O15").Value = 0.3484294080560655
$ws.Range("P15").Value = 0.3484294080560656
$ws.Range("Q15").Value = 572.2164711922169
$ws.Range("R15").Value = 5149.948240729952
$ws.Range("S15").Value = 0.04660376078732741
$ws.Range("T15").Value = 0.04660376078732743

$ws.Range("G16").Value = 14.14774133333333
$ws.Range("H16").Value = 42.443224
$ws.Range("I16").Value = 0.1337538098386587
$ws.Range("J16").Value = 0.1337538098386588
$ws.Range("M16").Value = 12.761795
$ws.Range("N16").Value = 38.28538500000001
$ws.Range("O16").Value = 0.1099393900775594
$ws.Range("P16").Value = 0.1099393900775594
$ws.Range("Q16").Value = 180.5505746090267
$ws.Range("R16").Value = 1624.95517148124
$ws.Range("S16").Value = 0.014704812274212
$ws.Range("T16").Value = 0.01470481227421201

$ws.Range("G17").Value = 14.14774133333333
$ws.Range("H17").Value = 42.443224
$ws.Range("I17").Value = 0.1337538098386587
$ws.Range("J17").Value = 0.1337538098386588
$ws.Range("M17").Value = 16.966758
$ws.Range("N17").Value = 50.900274
$ws.Range("O17").Value = 0.1461640017030168
$ws.Range("P17").Value = 0.1461640017030168
$ws.Range("Q17").Value = 240.041303449264
$ws.Range("R17").Value = 2160.371731043376
$ws.Range("S17").Value = 0.0195499920890427
$ws.Range("T17").Value = 0.01954999208904271
